$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Wine to discover (todo)")

# Rename the sheet: drop the "(todo)" suffix now that the deck is filled in.
$ws.Name = "Wine to discover"

# Remove the stray scratch cells that used to sit next to row 2.
$ws.Range("J2:M2").ClearContents()

# Fill in the remaining flashcards (rows 33-59), replacing the lone
# "000 Add (Red)" placeholder that used to live at A36.
$arr = New-Object 'object[,]' 27,2
$arr[0,0] = 'Saumur-Champigny Clos Rougeard-Foucault & fils'
$arr[0,1] = 'Cabernet franc, Loire, France'
$arr[1,0] = 'Saumur-Champigny Thierry Germain'
$arr[1,1] = 'Cabernet franc, Loire, France'
$arr[2,0] = 'Chinon Bernard'
$arr[2,1] = 'Baudry, Cabernet franc, Loire, France'
$arr[3,0] = 'Saint-Julien Léoville Barton'
$arr[3,1] = 'Cabernet sauvignon, merlot, Médoc, France'
$arr[4,0] = 'Savigny-Lès-Beaune 1er Cru Domaine de la Vougeraie'
$arr[4,1] = 'Pinot noir, Bourgogne, France'
$arr[5,0] = 'Pommard Grands Épenots Domaine de Courcel'
$arr[5,1] = 'Pinot noir, Bourgogne, France'
$arr[6,0] = 'Crozes Hermitage Emmanuel Darnaud'
$arr[6,1] = 'Syrah, Côtes du Rhône nord, France'
$arr[7,0] = 'Côte-Rôtie Domaine Jamet'
$arr[7,1] = 'Syrah, Côtes du Rhône nord, France'
$arr[8,0] = 'Châteauneuf-du-Pape Château Rayas'
$arr[8,1] = 'Grenache, Côtes du Rhône sud, France'
$arr[9,0] = 'Châteauneuf-du-Pape Henri Bonneau'
$arr[9,1] = 'Grenache, Côtes du Rhône sud, France'
$arr[10,0] = 'Bandol Château de Pibarnon'
$arr[10,1] = 'Mourvèdre, Provence, France'
$arr[11,0] = 'Barolo Giovanni Rosso'
$arr[11,1] = 'Nebbiolo, Piémont, Italy'
$arr[12,0] = 'Barolo Roberto Voerzio'
$arr[12,1] = 'Nebbiolo, Piémont, Italy'
$arr[13,0] = 'Chianti Classico La Massa'
$arr[13,1] = 'Sangiovese, Toscane, Italy'
$arr[14,0] = 'Brunello di Montalcino riserva Case Basse di Gianfranco Soldera'
$arr[14,1] = 'Sangiovese, Toscane, Italy'
$arr[15,0] = 'Alvaro Palacios'
$arr[15,1] = 'Carignan, Priorat, Spain'
$arr[16,0] = 'Vega Sicilia'
$arr[16,1] = 'Tempranillo, Ribera del Duero, Spain'
$arr[17,0] = 'Don Melchior'
$arr[17,1] = 'Cabernet sauvignon, Chile'
$arr[18,0] = 'Cheval des Andes'
$arr[18,1] = 'Malbec, Mendoza, Argentina'
$arr[19,0] = 'Henschke Hill of Grace'
$arr[19,1] = 'Syrah, Australia'
$arr[20,0] = 'Cordoba'
$arr[20,1] = 'Pinotage, Stellenbosch, South Africa'
$arr[21,0] = 'Beaux Frères'
$arr[21,1] = 'Pinot noir, Oregon, USA'
$arr[22,0] = 'Sloan Wines'
$arr[22,1] = 'Cabernet sauvignon, Napa Valley Rutherford, California'
$arr[23,0] = 'Screaming Eagle'
$arr[23,1] = 'Cabernet sauvignon, Napa Valley Rutherford, California'
$arr[24,0] = 'Dunn Vineyards Napa Valley'
$arr[24,1] = 'Cabernet sauvignon, Napa Valley Rutherford, California'
$arr[25,0] = 'Colgin Cellars Tychson Hill Vineyard'
$arr[25,1] = 'Cabernet sauvignon, Napa Valley Rutherford, California'
$arr[26,0] = 'Abreu'
$arr[26,1] = 'Cabernet sauvignon, Napa Valley Rutherford, California'
$ws.Range("A33:B59").Value = $arr

# Match the author's final scroll position/selection on this sheet.
$ws.Range("I2:M2").Select()

Write-Output "Wine to discover sheet rebuilt"
